$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1538
$ws.Range("F3").Value = 3295
$ws.Range("F5").Value = 723
$ws.Range("F6").Value = 2323
$ws.Range("F7").Value = 498
$ws.Range("F10").Value = 139
$ws.Range("F11").Value = 353
$ws.Range("F12").Value = 1100
$ws.Range("F13").Value = 454
$ws.Range("F15").Value = 89
$ws.Range("F16").Value = 255
$ws.Range("F17").Value = 4755
$ws.Range("F19").Value = 1345
$ws.Range("F20").Value = 3526
$ws.Range("F22").Value = 132
$ws.Range("F23").Value = 193
$ws.Range("F24").Value = 3758
$ws.Range("F25").Value = 5138
$ws.Range("F28").Value = 565
$ws.Range("F29").Value = 3303
$ws.Range("F30").Value = 379
$ws.Range("F35").Value = 1204
$ws.Range("F37").Value = 17
$ws.Range("F38").Value = 1428
$ws.Range("F40").Value = 1395
$ws.Range("F41").Value = 891
$ws.Range("F42").Value = 868
$ws.Range("F45").Value = 351
$ws.Range("F47").Value = 169
$ws.Range("F49").Value = 3742

# Sheet 2: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 3
$ws.Range("F6").Value = 1020
$ws.Range("F15").Value = 7

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 2311

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 2311
$ws.Range("F3").Value = 1538
$ws.Range("F4").Value = 3295
$ws.Range("F6").Value = 723
$ws.Range("F7").Value = 2323
$ws.Range("F8").Value = 498
$ws.Range("F11").Value = 1020
$ws.Range("F12").Value = 139
$ws.Range("F13").Value = 353
$ws.Range("F14").Value = 1100
$ws.Range("F15").Value = 454
$ws.Range("F17").Value = 89
$ws.Range("F18").Value = 255
$ws.Range("F19").Value = 4755
$ws.Range("F20").Value = 1345
$ws.Range("F21").Value = 3526
$ws.Range("F22").Value = 3758
$ws.Range("F23").Value = 5138
$ws.Range("F25").Value = 565
$ws.Range("F26").Value = 3303
$ws.Range("F27").Value = 379
$ws.Range("F32").Value = 1204
$ws.Range("F34").Value = 17
$ws.Range("F35").Value = 1428
$ws.Range("F36").Value = 1395
$ws.Range("F37").Value = 891
$ws.Range("F43").Value = 351
$ws.Range("F46").Value = 169
$ws.Range("F49").Value = 3742
